$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.575.98"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "2.333.28"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.536"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("D9").Value = "2.367.74"
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.103"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.153"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.353"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "2.757.18"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "57.649.22"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "2.346.97"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.98%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("E25").Value = "  +2.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.993"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "174.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.99%  "
$ws.Range("E30").Value = "  +5.00%  "
$ws.Range("D31").Value = "0.0₃0739"
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("B34").Value = "SuiNetwork"
$ws.Range("C34").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +13.13%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.989"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "39.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "150.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.378"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "285.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0929"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0506"
$ws.Range("D46").Style = "Normal"
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.98%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.563"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("E50").Value = "  +2.71%  "
$ws.Range("E51").Value = "  +6.90%  "
